# Generate Report for Handoff
# The "2c29110a-a729-4e33-bcb9-89a01b98271b.md" file (row 3 in every sheet)
# has been re-handed-off: its status flips from "Handed back: in sync with
# en-US" to "Ready for handoff", its timestamps advance, and the localized
# sheets record that the handback file on disk is stale (new Error Detail
# text + widened column P).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: summary status/date columns for the 2c29110a... row
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-23 14:58:36"

# ---------------------------------------------------------------------
# zh-cn sheet: status, handoff datetime and new error detail for row 3
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("H3").Value = "2016-08-23 14:58:31"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0bad0842abb4c726cbac78ab460f1c358f4ca0d7/e2e/2c29110a-a729-4e33-bcb9-89a01b98271b.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9759ce726f2bfc22247efff27ef3fe97478e9b57/e2e/2c29110a-a729-4e33-bcb9-89a01b98271b.md."
# Widen the Error Detail column so the long message is readable
# (input value compensates for Excel's column-width padding so the
# persisted OOXML width lands on exactly 40).
$zhcn.Columns.Item(16).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------
# de-de sheet: status, handoff datetime and new error detail for row 3
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("H3").Value = "2016-08-23 14:58:36"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0bad0842abb4c726cbac78ab460f1c358f4ca0d7/e2e/2c29110a-a729-4e33-bcb9-89a01b98271b.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9759ce726f2bfc22247efff27ef3fe97478e9b57/e2e/2c29110a-a729-4e33-bcb9-89a01b98271b.md."
$dede.Columns.Item(16).ColumnWidth = 39.166666666666664
